$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# "SUCCESS:CT,0:{ changePercent, changeTolerance,xHighMax,xLowMax,yHighMax,yLowMax}"
# becomes three runs:
#   "SUCCESS:CT,0:{changePercent,"  (sz/szCs 18, unchanged formatting)
#   " "                              (a plain space, different/default formatting)
#   "xHighChangeTolerance,xLowChangeTolerance, yHighChangeTolerance,yLowChangeTolerance}" (sz/szCs 18)
$findRng = $d.Content
$findRng.Find.Execute(
    "SUCCESS:CT,0:{ changePercent, changeTolerance,xHighMax,xLowMax,yHighMax,yLowMax}",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($findRng.Find.Found) {
    $cellStart = $findRng.Start
    $newText = "SUCCESS:CT,0:{changePercent, xHighChangeTolerance,xLowChangeTolerance, yHighChangeTolerance,yLowChangeTolerance}"
    $findRng.Text = $newText

    # locate the single space that sits between "changePercent," and "xHighChangeTolerance"
    $beforeSpace = "SUCCESS:CT,0:{changePercent,"
    $spaceStart = $cellStart + $beforeSpace.Length
    $spaceRange = $d.Range($spaceStart, $spaceStart + 1)

    # give that space run its own (default/document) formatting, splitting the
    # single run into three runs, matching what a manual retype of the space
    # would produce
    $spaceRange.Font.Size = 11
}

# --- Change 2 ---------------------------------------------------------
# "Get change tolerance value based..." -> "Get change tolerance values based..."
$d.Content.Find.Execute(
    "Get change tolerance value based on max value of FSRs and change tolerance percentage",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Get change tolerance values based on max value of FSRs and change tolerance percentage",
    2)
